# Update cryptos list values per latest data pull
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.191.68"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.636.87"
$ws.Range("E3").Value = "  -0.89%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.04"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.09"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.31%  "

$ws.Range("E9").Value = "  -5.02%  "

$ws.Range("E10").Value = "  +0.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.336"
$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.103.91"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.226.24"
$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.70"
$ws.Range("E15").Value = "  -2.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.652.43"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "343.16"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("E19").Value = "  +0.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.55"
$ws.Range("E20").Value = "  +1.51%  "

$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.48"
$ws.Range("E23").Value = "  +4.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.414"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.766.35"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0793"
$ws.Range("E29").Value = "  -1.40%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.31"
$ws.Range("E31").Value = "  -5.08%  "

$ws.Range("E32").Value = "  +0.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.95"
$ws.Range("E33").Value = "  +0.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.92"
$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.13"
$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("E36").Value = "  -2.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.851"
$ws.Range("E37").Value = "  -5.05%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.30"
$ws.Range("E38").Value = "  -1.19%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.843"
$ws.Range("E39").Value = "  -4.59%  "

$ws.Range("E40").Value = "  -2.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.62"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0978"
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.598"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "268.95"
$ws.Range("E45").Value = "  -2.63%  "

$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.19"
$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("E48").Value = "  -1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.035.93"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.74"
$ws.Range("E51").Value = "  -1.73%  "
